# "MobileNetV3 (large) was selected."
# Fill in the Accuracy / Loss benchmark numbers for every model (the
# MobileNetV3(large) row, which previously had no figures at all, now gets
# its Accuracy/Loss too) and tidy up the layout that goes with it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The whole table had picked up stray explicit cell styles (left-aligned /
# "0.000_ " number format) along the way; the refreshed data drops those and
# falls back to the sheet's default style, and the blank corner cell (A1)
# goes back to being completely empty.
$ws.Range("A1:E11").ClearFormats()
$ws.Range("A1").ClearContents()

# Column B (Accuracy) is a touch narrower now that it holds real numbers.
$ws.Columns.Item(2).ColumnWidth = 4.86

# Accuracy (B) / Loss (C) for every model, keyed by row number.
$results = @{
  2  = @(0.73988103866577104, 0.72463768720626798)   # MobileNet(alpha=0.25)
  3  = @(0.45119205117225603, 0.84601449966430597)   # MobileNetV2(alpha=0.35)
  4  = @(0.33470913767814597, 0.88949275016784601)   # MobileNetV2(alpha=0.50)
  5  = @(0.524896681308746,   0.80344200134277299)   # MobileNet(alpha=0.50)
  6  = @(0.27091634273528997, 0.90670287609100297)   # MobileNetV3(small)
  7  = @(0.25743499398231501, 0.90942031145095803)   # MobileNetV2(alpha=0.75)
  8  = @(0.47506558895111001, 0.82336956262588501)   # MobileNet(alpha=0.75)
  9  = @(0.28164747357368403, 0.90307968854904097)   # MobileNetV2(alpha=1.0)
  10 = @(0.24787789583206099, 0.92028987407684304)   # MobileNetV3(large)  <- newly selected
  11 = @(0.42793497443199102, 0.85235506296157804)   # MobileNet(alpha=1.0)
}

foreach ($row in $results.Keys) {
  $pair = $results[$row]
  $ws.Cells.Item($row, 2).Value = $pair[0]
  $ws.Cells.Item($row, 3).Value = $pair[1]
}

# New, two-decimal number format for the freshly filled Accuracy/Loss cells.
$ws.Range("B2:C11").NumberFormat = "0.00_ "
